# Common: A lot of small fixes/improvements
# Adds new translation rows (lab.* menu labels) to the "Import" sheet,
# mirroring the formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")
$ws.Activate()

$data = @(
    @("lab.vape.menu", "Vapování"),
    @("lab.liquid.menu", "Liquidy"),
    @("lab.atomizer.menu", "Atomizéry"),
    @("lab.mod.menu", "Mody"),
    @("lab.cell.menu", "Články"),
    @("lab.vendor.menu", "Výrobci"),
    @("lab.cotton.menu", "Vaty"),
    @("lab.wire.menu", "Dráty"),
    @("lab.coil.menu", "Spirálky")
)

$lastRow = 272
$startRow = $lastRow + 1
$endRow = $startRow + $data.Count - 1

# Copy formatting from the last existing data row down across the whole new
# block so the new rows match the sheet's style (font, wrap, column styles).
$ws.Range("A$lastRow`:C$lastRow").Copy() | Out-Null
$ws.Range("A$startRow`:C$endRow").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
}

# Reflect the final cursor/selection position from the edit session.
$ws.Range("B277").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 269
$win.ScrollColumn = 1
